$p = $ppt.ActivePresentation

# Insert a new slide right after slide 5 ("I THINK RN") / before slide 6 ("ISSUES"),
# i.e. at new index 6, using the same "Title and Content" layout as the ISSUES slide.
$layout = $p.Slides.Item(6).CustomLayout
$newSlide = $p.Slides.AddSlide(6, $layout)

# Title placeholder stays empty (just an empty paragraph) - nothing to type.

# Content placeholder (shape 2 on the new slide) gets the note about reverting the
# zombie animation.
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "If I need to revert the z" + [char]92 + "ombie animiaton then delete links to attack and "

# Add a textbox with the C# zombie-attack script snippet.
$code = "using UnityEngine;`r" + `
    "`r" + `
    "public class ZombieAttack : MonoBehaviour`r" + `
    "{`r" + `
    "    public int damage = 10; // Damage per hit`r" + `
    "    public float attackRate = 1.5f; // Attack every X seconds`r" + `
    "`r" + `
    "    private float nextAttackTime = 0f;`r" + `
    "`r" + `
    "    private void OnCollisionStay(Collision collision)`r" + `
    "    {`r" + `
    "        if (collision.gameObject.CompareTag(" + [char]34 + "Player" + [char]34 + "))`r" + `
    "        {`r" + `
    "            if (Time.time >= nextAttackTime)`r" + `
    "            {`r" + `
    "                HealthManager playerHealth = collision.gameObject.GetComponent<HealthManager>();`r" + `
    "                DamageOverlay screenEffect = collision.gameObject.GetComponent<DamageOverlay>(); // Get the red pulse effect`r" + `
    "`r" + `
    "                if (playerHealth != null)`r" + `
    "                {`r" + `
    "                    playerHealth.TakeDamage(damage);`r" + `
    "                    Debug.Log(" + [char]34 + "Zombie dealt " + [char]34 + " + damage + " + [char]34 + " damage to player." + [char]34 + ");`r" + `
    "                }`r" + `
    "`r" + `
    "                if (screenEffect != null)`r" + `
    "                {`r" + `
    "                    screenEffect.ShowDamageEffect(); // Trigger red screen pulse`r" + `
    "                }`r" + `
    "`r" + `
    "                nextAttackTime = Time.time + attackRate;`r" + `
    "            }`r" + `
    "        }`r" + `
    "    }`r" + `
    "}"

$tb = $newSlide.Shapes.AddTextbox(1, 368.64, 198.895, 484.56, 336.858)
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0
$tb.TextFrame.TextRange.Text = $code
$tb.TextFrame.TextRange.Font.Size = 8
$tb.Left = 368.64
$tb.Top = 198.895
$tb.Width = 484.56
$tb.Height = 336.858
